$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header column in H1, reusing the same formatting as the
# other header cells (e.g. G1 "sum") by copying formats over.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Add the numeric values for the new Save column
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
